# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback Datetime" (column H) for the first file entry
# (row 2, 2b08ecc8-b4ee-478b-a966-1af50d4e1bca.md) on both the zh-cn and
# de-de language report sheets, reflecting the latest handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 04:50:36"
$wsZhCn.Range("H2").Value = "2016-03-24 04:51:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 04:50:40"
$wsDeDe.Range("H2").Value = "2016-03-24 04:51:08"
